# Re-run of the Data Preprocessing / Cleaning script.
# This refreshes the "timestamp" column (N) for every data row with the
# new run's timestamps, and updates a handful of re-fetched / recomputed
# fields (subregion, capital, population, area, population_density) that
# changed between the two script runs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Field-level corrections produced by the re-run (subregion/capital
#    lookups that resolved differently, and population/area/density
#    recomputation for a few countries).
# ---------------------------------------------------------------------
$ws.Range("F4").Value  = "Western Europe"          # Switzerland subregion
$ws.Range("J5").Value  = "Unknown"                  # Sierra Leone capital

$ws.Range("H12").Value = 163610                     # Tunisia area
$ws.Range("O12").Value = 72.23652588472588          # Tunisia population_density

$ws.Range("J14").Value = "Unknown"                  # Benin capital
$ws.Range("F24").Value = "Unknown"                  # France subregion

$ws.Range("G29").Value = 2077132                    # North Macedonia population
$ws.Range("O29").Value = 80.78139462528682          # North Macedonia population_density

$ws.Range("G68").Value = 1026                       # Antarctica population
$ws.Range("O68").Value = 0.00007328571428571429     # Antarctica population_density

$ws.Range("H76").Value = 897831.3548002663          # Pakistan area
$ws.Range("O76").Value = 246.0287556443606          # Pakistan population_density

$ws.Range("H113").Value = 50180.73081625393         # Bosnia and Herzegovina area
$ws.Range("O113").Value = 65.3799764697193          # Bosnia and Herzegovina population_density

$ws.Range("J128").Value = "[""Yaren""]"             # Nauru capital

$ws.Range("G130").Value = 517714                    # Western Sahara population
$ws.Range("O130").Value = 1.946293233082707         # Western Sahara population_density

$ws.Range("H197").Value = 1240192                   # Mali area
$ws.Range("O197").Value = 16.32878941325214         # Mali population_density

$ws.Range("F219").Value = "Unknown"                 # Kuwait subregion
$ws.Range("J220").Value = "[""Mal" + [char]0x00e9 + """]"   # Maldives capital ("Malé")

$ws.Range("G224").Value = 212559409                 # Brazil population
$ws.Range("O224").Value = 24.960688684883           # Brazil population_density

$ws.Range("F229").Value = "South America"           # Venezuela subregion

# ---------------------------------------------------------------------
# 2) Refresh the "timestamp" column (N) for every data row (2-251) with
#    the new run's batch timestamps.
# ---------------------------------------------------------------------
$tsGroups = @(
    @{ start=2;   end=8;   ts="2025-03-31T00:39:25.402721" },
    @{ start=9;   end=16;  ts="2025-03-31T00:39:25.418344" },
    @{ start=17;  end=25;  ts="2025-03-31T00:39:25.434479" },
    @{ start=26;  end=31;  ts="2025-03-31T00:39:25.450149" },
    @{ start=32;  end=39;  ts="2025-03-31T00:39:25.465776" },
    @{ start=40;  end=48;  ts="2025-03-31T00:39:25.481418" },
    @{ start=49;  end=56;  ts="2025-03-31T00:39:25.497026" },
    @{ start=57;  end=65;  ts="2025-03-31T00:39:25.512652" },
    @{ start=66;  end=73;  ts="2025-03-31T00:39:25.528276" },
    @{ start=74;  end=82;  ts="2025-03-31T00:39:25.543906" },
    @{ start=83;  end=90;  ts="2025-03-31T00:39:25.559531" },
    @{ start=91;  end=96;  ts="2025-03-31T00:39:25.575148" },
    @{ start=97;  end=104; ts="2025-03-31T00:39:25.590776" },
    @{ start=105; end=113; ts="2025-03-31T00:39:25.606400" },
    @{ start=114; end=122; ts="2025-03-31T00:39:25.622033" },
    @{ start=123; end=131; ts="2025-03-31T00:39:25.637651" },
    @{ start=132; end=139; ts="2025-03-31T00:39:25.653274" },
    @{ start=140; end=147; ts="2025-03-31T00:39:25.668899" },
    @{ start=148; end=155; ts="2025-03-31T00:39:25.684525" },
    @{ start=156; end=162; ts="2025-03-31T00:39:25.700149" },
    @{ start=163; end=171; ts="2025-03-31T00:39:25.715777" },
    @{ start=172; end=180; ts="2025-03-31T00:39:25.731399" },
    @{ start=181; end=188; ts="2025-03-31T00:39:25.747086" },
    @{ start=189; end=197; ts="2025-03-31T00:39:25.762651" },
    @{ start=198; end=205; ts="2025-03-31T00:39:25.778274" },
    @{ start=206; end=213; ts="2025-03-31T00:39:25.793899" },
    @{ start=214; end=221; ts="2025-03-31T00:39:25.809529" },
    @{ start=222; end=230; ts="2025-03-31T00:39:25.825157" },
    @{ start=231; end=239; ts="2025-03-31T00:39:25.840774" },
    @{ start=240; end=247; ts="2025-03-31T00:39:25.856399" },
    @{ start=248; end=251; ts="2025-03-31T00:39:25.872025" }
)

foreach ($g in $tsGroups) {
    for ($r = $g.start; $r -le $g.end; $r++) {
        $ws.Cells.Item($r, 14).Value = $g.ts
    }
}
